$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 270
$ws.Range("I6").Value = 292.5
$ws.Range("K6").Value = 877.5
$ws.Range("M6").Value = -765.5
$ws.Range("H33").Value = 681.6875
$ws.Range("I33").Value = 736.2143
$ws.Range("K33").Value = 736.2143
$ws.Range("M33").Value = -507.2143
$ws.Range("H55").Value = 635.3913
$ws.Range("I55").Value = 631.75
$ws.Range("K55").Value = 631.75
$ws.Range("M55").Value = -417.75
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
$ws.Range("H100").Value = 7369.381
$ws.Range("I100").Value = 6002.778
$ws.Range("K100").Value = 6002.778
$ws.Range("M100").Value = -5461.778
$ws.Range("H131").Value = 6012.364
$ws.Range("I131").Value = 3613.6
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 10840.8
$ws.Range("L131").Value = 90000
$ws.Range("M131").Value = -5800.799999999999
$ws.Range("N131").Value = -100080
$ws.Range("H135").Value = 3189.9412
$ws.Range("I135").Value = 3479.6
$ws.Range("J135").Value = 1017.5
$ws.Range("K135").Value = 31316.4
$ws.Range("L135").Value = 9157.5
$ws.Range("M135").Value = -28781.4
$ws.Range("N135").Value = -14227.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5572.339
$ws.Range("I32").Value = 4944.293
$ws.Range("K32").Value = 4944.293
$ws.Range("M32").Value = -4657.293

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4224.933
$ws.Range("I20").Value = 3126
$ws.Range("J20").Value = 7247
$ws.Range("K20").Value = 3126
$ws.Range("L20").Value = 7247
$ws.Range("M20").Value = -2879
$ws.Range("N20").Value = -7741
$ws.Range("H94").Value = 6669864.5
$ws.Range("I94").Value = 2405.923
$ws.Range("J94").Value = 13892945
$ws.Range("K94").Value = 2405.923
$ws.Range("L94").Value = 13892945
$ws.Range("M94").Value = -1954.923
$ws.Range("N94").Value = -13893847
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 68750
$ws.Range("J117").Value = 81666.664
$ws.Range("L117").Value = 81666.664
$ws.Range("N117").Value = -90844.664
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H119").Value = 44991
$ws.Range("J119").Value = 44991
$ws.Range("L119").Value = 44991
$ws.Range("N119").Value = -54667
$ws.Range("H120").Value = 641666.3
$ws.Range("J120").Value = 641666.3
$ws.Range("L120").Value = 641666.3
$ws.Range("N120").Value = -648924.3
$ws.Range("H132").Value = 5244.75
$ws.Range("I132").Value = 4326.3335
$ws.Range("K132").Value = 12979.0005
$ws.Range("M132").Value = -10449.0005
$ws.Range("H139").Value = 89982.5
$ws.Range("J139").Value = 89982.5
$ws.Range("L139").Value = 89982.5
$ws.Range("N139").Value = -100262.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.77778000000001
$ws.Range("I2").Value = 11.6
$ws.Range("K2").Value = 69.59999999999999
$ws.Range("M2").Value = 43.40000000000001
$ws.Range("H5").Value = 3799.2334
$ws.Range("I5").Value = 3144.2632
$ws.Range("J5").Value = 4930.5454
$ws.Range("K5").Value = 9432.7896
$ws.Range("L5").Value = 14791.6362
$ws.Range("M5").Value = -9320.7896
$ws.Range("N5").Value = -15015.6362
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 3
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 9
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 422
$ws.Range("N47").ClearContents()
$ws.Range("H50").Value = 2525.5
$ws.Range("J50").Value = 3699.75
$ws.Range("L50").Value = 11099.25
$ws.Range("N50").Value = -12061.25
$ws.Range("H53").Value = 2525.5
$ws.Range("J53").Value = 3699.75
$ws.Range("L53").Value = 11099.25
$ws.Range("N53").Value = -12061.25
$ws.Range("H68").Value = 817.36365
$ws.Range("I68").Value = 761.375
$ws.Range("J68").Value = 966.6667
$ws.Range("K68").Value = 2284.125
$ws.Range("L68").Value = 2900.0001
$ws.Range("M68").Value = -1473.125
$ws.Range("N68").Value = -4522.0001
$ws.Range("H71").Value = 817.36365
$ws.Range("I71").Value = 761.375
$ws.Range("J71").Value = 966.6667
$ws.Range("K71").Value = 6852.375
$ws.Range("L71").Value = 8700.0003
$ws.Range("M71").Value = -2796.375
$ws.Range("N71").Value = -16812.0003
$ws.Range("H92").Value = 665
$ws.Range("I92").Value = 665
$ws.Range("K92").Value = 1995
$ws.Range("M92").Value = -747
$ws.Range("H104").Value = 3471.2
$ws.Range("J104").Value = 8663.333000000001
$ws.Range("L104").Value = 25989.999
$ws.Range("N104").Value = -31231.999
$ws.Range("H132").Value = 1363
$ws.Range("I132").Value = 1264.6666
$ws.Range("J132").Value = 1399.875
$ws.Range("K132").Value = 11381.9994
$ws.Range("L132").Value = 12598.875
$ws.Range("M132").Value = -8851.999400000001
$ws.Range("N132").Value = -17658.875
$ws.Range("H134").Value = 6176.5557
$ws.Range("I134").Value = 2227
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 6681
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -1611
$ws.Range("N134").Value = -70140
$ws.Range("H135").Value = 3799.2334
$ws.Range("I135").Value = 3144.2632
$ws.Range("J135").Value = 4930.5454
$ws.Range("K135").Value = 28298.3688
$ws.Range("L135").Value = 44374.9086
$ws.Range("M135").Value = -25763.3688
$ws.Range("N135").Value = -49444.9086
$ws.Range("H137").Value = 1671.8
$ws.Range("I137").Value = 1671.8
$ws.Range("K137").Value = 5015.4
$ws.Range("M137").Value = 84.60000000000036
$ws.Range("H140").Value = 2396.7058
$ws.Range("I140").Value = 2396.7058
$ws.Range("K140").Value = 7190.117400000001
$ws.Range("M140").Value = -2010.117400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1283.5
$ws.Range("I132").Value = 1337.9375
$ws.Range("J132").Value = 848
$ws.Range("K132").Value = 4013.8125
$ws.Range("L132").Value = 2544
$ws.Range("M132").Value = -1483.8125
$ws.Range("N132").Value = -7604

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4284.25
$ws.Range("I82").Value = 2466.4375
$ws.Range("J82").Value = 7919.875
$ws.Range("K82").Value = 2466.4375
$ws.Range("L82").Value = 7919.875
$ws.Range("M82").Value = -2105.4375
$ws.Range("N82").Value = -8641.875
$ws.Range("H85").Value = 4284.25
$ws.Range("I85").Value = 2466.4375
$ws.Range("J85").Value = 7919.875
$ws.Range("K85").Value = 2466.4375
$ws.Range("L85").Value = 7919.875
$ws.Range("M85").Value = -1218.4375
$ws.Range("N85").Value = -10415.875
$ws.Range("H132").Value = 3703.0667
$ws.Range("I132").Value = 1963.0588
$ws.Range("K132").Value = 5889.1764
$ws.Range("M132").Value = -3359.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14650.182
$ws.Range("J41").Value = 13942.2
$ws.Range("L41").Value = 13942.2
$ws.Range("N41").Value = -14722.2
$ws.Range("H81").Value = 11909590
$ws.Range("I81").Value = 5740.087
$ws.Range("K81").Value = 11480.174
$ws.Range("M81").Value = -10419.174
$ws.Range("H84").Value = 11909590
$ws.Range("I84").Value = 5740.087
$ws.Range("K84").Value = 57400.87
$ws.Range("M84").Value = -52096.87
$ws.Range("H96").Value = 1925.4445
$ws.Range("I96").Value = 1910.4
$ws.Range("K96").Value = 1910.4
$ws.Range("M96").Value = -537.4000000000001
$ws.Range("H100").Value = 568.6667
$ws.Range("I100").Value = 346.33334
$ws.Range("J100").Value = 1013.3333
$ws.Range("K100").Value = 692.66668
$ws.Range("L100").Value = 2026.6666
$ws.Range("M100").Value = -151.66668
$ws.Range("N100").Value = -3108.6666
$ws.Range("H132").Value = 2963.5
$ws.Range("I132").Value = 1972.625
$ws.Range("K132").Value = 5917.875
$ws.Range("M132").Value = -3387.875
$ws.Range("H141").Value = 93846
$ws.Range("J141").Value = 93846
$ws.Range("L141").Value = 93846
$ws.Range("N141").Value = -104206
